# Apply "Fixed R, Tableau and Physics" edits to the "skill" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("skill")

# Fix Physics: row 38 (Skill_ID 36) category changes from Hard -> Education
$ws.Range("D38").Value = "Education"

# Row 50 ("Physics"/"Major") was a stray duplicate of the Physics entry
# already fixed above. Shift the Skill/Category text of rows 51-56 up into
# rows 50-55 (the Skill_ID and JobTitle columns are untouched), then clear
# the now-unused row 56.
$ws.Range("B50").Value = "Economics"
$ws.Range("B51").Value = "Information System"
$ws.Range("B52").Value = "Quantitative Finance"
$ws.Range("B53").Value = "Biostatistics"
$ws.Range("B54").Value = "Bioinformatics"
$ws.Range("B55").Value = "MBA"
$ws.Range("D55").Value = "Education"

$ws.Rows.Item(56).Delete()

# Select the new last-used view state as recorded after the edit.
$ws.Range("A21").Select()
$excel.ActiveWindow.ScrollRow = 21
$ws.Range("D39").Select()
